# Assign student to exam fixed to a specific school:
#  - rename the existing "Nilai Siswa" sheet to "ujian 1"
#  - drop the "Nilai 1"/"Nilai 2" columns (G:H) that are no longer needed
#  - add a second sheet "ujian 2" with the same (trimmed) layout/data

$wb = $excel.ActiveWorkbook

# 1) Rename the first sheet and strip the trailing Nilai 1 / Nilai 2 columns.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ujian 1"
$ws1.Columns.Item(7).Delete()
$ws1.Columns.Item(7).Delete()

# 2) Add the second exam sheet right after the first one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ujian 2"

# 3) Populate it with the same headers/data (values + styles) as sheet one.
$ws1.Range("A1:F2").Copy($ws2.Range("A1"))
